$d = $word.ActiveDocument

# The document has one section whose "first page" header/footer are
# distinct from the "default" (primary) header/footer:
#   Headers.Item(1) / Footers.Item(1) -> wdHeaderFooterPrimary (default)
#   Headers.Item(2) / Footers.Item(2) -> wdHeaderFooterFirstPage (first)
#
# Mapping of the inline pictures in this document:
#   Footers.Item(1)  (primary / footer2.xml) -> Pearson Edexcel logo, id=2
#   Footers.Item(2)  (first page / footer1.xml) -> Pearson Edexcel logo, id=3
#   Headers.Item(2)  (first page / header1.xml) -> BTEC logo, id=1

$section = $d.Sections.Item(1)

# Pearson Edexcel logo in the primary (default) footer: image2.png -> image1.png
$footerPrimary = $section.Footers.Item(1)
$pearsonShape1 = $footerPrimary.Range.InlineShapes.Item(1)
$pearsonShape1.Name = "image1.png"

# Pearson Edexcel logo in the first-page footer: image2.png -> image1.png
$footerFirst = $section.Footers.Item(2)
$pearsonShape2 = $footerFirst.Range.InlineShapes.Item(1)
$pearsonShape2.Name = "image1.png"

# BTEC logo in the first-page header: image1.jpg -> image2.jpg
$headerFirst = $section.Headers.Item(2)
$btecShape = $headerFirst.Range.InlineShapes.Item(1)
$btecShape.Name = "image2.jpg"
